# Scrum board update: fold in everyone's contributions under the "Jay" /
# login-feature swim-lane (columns E:F), expanding it from 5 rows to 12
# rows, and push the rest of the board (rows 7-38) down by 7 rows to make
# room.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# 1. Make room: insert 7 blank rows above the old row 7 ("As Jim..." story)
#    so everything that used to live at rows 7-38 now lives at rows 14-45.
# ---------------------------------------------------------------------
$ws.Rows("7:13").Insert()

# ---------------------------------------------------------------------
# 2. Rewrite the E:F "Done" swim-lane (rows 2-13) with the full task list
#    and who did it.
# ---------------------------------------------------------------------

# Column E: task description text for each row.
$tasks = @{
    2  = "Create a database of users"
    3  = "Implement a way to hash their user information"
    4  = "Store all user information"
    5  = "Create a login page"
    6  = "Created a splash screen for the application"
    7  = "Created an onClickListener for the login page"
    8  = "Created edit texts and button for the login page"
    9  = "Created a server application for the application to connect to"
    10 = "Added a command system in the server to check for command received by client"
    11 = "Implemented login command for the server (read user pass and return verification to client)"
    12 = "Implemented the login activity's onClickListener to talk to the server"
    13 = "Login activity receives whether or not a login was successful"
}

# Column F: who worked on it.
$owners = @{
    2  = "Richard"
    3  = "Richard"
    4  = "Charlie"
    5  = "Jay"
    6  = "Gregory"
    7  = "Jay"
    8  = "Bo Fang"
    9  = "Gregory"
    10 = "Gregory"
    11 = "Jay"
    12 = "Charlie"
    13 = "Bo Fang"
}

# Row heights for the rewritten block.
$heights = @{
    2  = 47.25
    3  = 31.5
    4  = 15.75
    5  = 15.75
    6  = 31.5
    7  = 31.5
    8  = 30
    9  = 30
    10 = 45
    11 = 45
    12 = 30
    13 = 30
}

# Per-person font colour (VBA BGR-encoded long, matches the workbook's
# existing colour coding convention).
$ownerColor = @{
    "Richard" = 255        # FFFF0000
    "Charlie" = 192         # FFC00000
    "Jay"     = 5287936     # FF00B050
    "Gregory" = 49407       # FFFFC000
    "Bo Fang" = 10498160    # FF7030A0
}

foreach ($r in 2..13) {
    $eCell = $ws.Cells.Item($r, 5)   # column E
    $fCell = $ws.Cells.Item($r, 6)   # column F

    $eCell.Value = $tasks[$r]
    $fCell.Value = $owners[$r]

    # Column E formatting: size-12 default-colour text, centered, wrapped.
    # Rows 2-7 are vertically centered too; rows 8-13 match the plainer
    # "wrap + horizontal-center only" look used for the newer entries.
    $eCell.Font.Name = "Calibri"
    $eCell.Font.Size = 12
    $eCell.Font.ColorIndex = -4105
    $eCell.HorizontalAlignment = -4108
    $eCell.WrapText = $true
    if ($r -le 7) {
        $eCell.VerticalAlignment = -4108
    } else {
        $eCell.Font.Size = 11
        $eCell.VerticalAlignment = -4142
    }

    # Column F formatting: size-11 text in the owner's colour, centered.
    $fCell.Font.Name = "Calibri"
    $fCell.Font.Size = 11
    $fCell.Font.Color = $ownerColor[$owners[$r]]
    $fCell.HorizontalAlignment = -4108
    $fCell.VerticalAlignment = -4108
    $fCell.WrapText = $false

    $ws.Rows($r).RowHeight = $heights[$r]
}

# Rows 6 and 7 never had a column-B cell in the new layout (unlike rows
# 8-13, which keep an empty styled B cell) - clear any stray formatting
# that row-insert may have copied down.
$ws.Range("B6:B7").ClearFormats()

# ---------------------------------------------------------------------
# 3. Selection follows where the editor was last working.
# ---------------------------------------------------------------------
$ws.Range("F3").Select()
